$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 14.950000000000003
$ws.Range("C2").Value = 11.450000000000001
$ws.Range("D2").Value = 11.799999999999999
$ws.Range("E2").Value = 13.25

# Row 3 data values
$ws.Range("B3").Value = 11.950000000000001
$ws.Range("C3").Value = 5.3000000000000007
$ws.Range("D3").Value = 9.75
$ws.Range("E3").Value = 11.15

# Update the selected range to match the new extent of updated data
$ws.Range("B1:E3").Select() | Out-Null
